$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 99 - this shifts the existing rows 99-117
# down to rows 101-119, matching the diff (dimension A1:R117 -> A1:R119).
$ws.Rows("99:100").Insert()

# New row 99: Agrícola del Norte S.A. de Arica - Caigua, Primera, week of 44785
$ws.Range("A99").Value = 1
$ws.Range("B99").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C99").Value = "Arica y Parinacota"
$ws.Range("D99").Value = 44785
$ws.Range("E99").Value = 15
$ws.Range("F99").Value = 100112036
$ws.Range("G99").Value = "Caigua"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 150
$ws.Range("K99").Value = 7000
$ws.Range("L99").Value = 8000
$ws.Range("M99").Value = 7500
$ws.Range("N99").Value = "`$/caja 20 kilos"
$ws.Range("O99").Value = "Región de Arica y Parinacota"
$ws.Range("P99").Value = 375
$ws.Range("Q99").Value = 20
$ws.Range("R99").Value = "Hortaliza"

# New row 100: Agrícola del Norte S.A. de Arica - Caigua, Segunda, week of 44785
$ws.Range("A100").Value = 1
$ws.Range("B100").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C100").Value = "Arica y Parinacota"
$ws.Range("D100").Value = 44785
$ws.Range("E100").Value = 15
$ws.Range("F100").Value = 100112036
$ws.Range("G100").Value = "Caigua"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Segunda"
$ws.Range("J100").Value = 160
$ws.Range("K100").Value = 6000
$ws.Range("L100").Value = 7000
$ws.Range("M100").Value = 6500
$ws.Range("N100").Value = "`$/caja 20 kilos"
$ws.Range("O100").Value = "Región de Arica y Parinacota"
$ws.Range("P100").Value = 325
$ws.Range("Q100").Value = 20
$ws.Range("R100").Value = "Hortaliza"

# Ensure date cells keep the workbook's date number format (style index 2,
# format "YYYY-MM-DD HH:MM:SS") matching the rest of column D.
$ws.Range("D99").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D100").NumberFormat = "YYYY-MM-DD HH:MM:SS"
